# Architecture general components.docx - apply diff:
#  1) Insert a blank paragraph (pPr/rPr lang eastAsia="en-NL", no list/style,
#     no run) right before the "backend server" bullet paragraph.
#  2) Insert the same kind of blank paragraph right before the "database"
#     bullet paragraph.
#  3) Split the "machine learning pipeline" bullet paragraph: insert two
#     blank paragraphs right before it - the first keeps the ListParagraph
#     style (but no numbering), the second has neither style nor numbering -
#     leaving the original bullet (with its numbering + text) intact after
#     them.
#  4) Add <w:lastRenderedPageBreak/> to the start of the "Postprocessing..."
#     run.
#  5) Remove <w:lastRenderedPageBreak/> from the start of the "Evaluation..."
#     run (it effectively moved from there to the Postprocessing run above).

$d = $word.ActiveDocument
$W_NS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaXml($rng) {
    # Round-trips a single paragraph Range through WordOpenXML and pulls out
    # just the <w:p>...</w:p> for that paragraph (the first one emitted -
    # Word always appends a filler paragraph + sectPr after it).
    $full = $rng.WordOpenXML
    if ($full -match '<w:body>(<w:p\b.*?</w:p>)') {
        return $matches[1]
    }
    return $null
}

function Insert-BlankBefore($findText, $fragmentsXml) {
    # Locates the paragraph that starts with $findText and prepends the
    # given blank-paragraph fragment(s) (array of XML strings) in front of
    # it, leaving the found paragraph's own content untouched.
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph starting with: $findText"
    }
    $para = $range.Paragraphs(1)
    $orig = Get-ParaXml($para.Range)
    $prefix = [string]::Join("", $fragmentsXml)
    $para.Range.InsertXML($prefix + $orig)
}

$blankNoStyle = "<w:p $W_NS><w:pPr><w:rPr><w:lang w:eastAsia=`"en-NL`"/></w:rPr></w:pPr></w:p>"
$blankListStyle = "<w:p $W_NS><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:rPr><w:lang w:eastAsia=`"en-NL`"/></w:rPr></w:pPr></w:p>"

# 1) before "The backend server would handle the requests from the client..."
Insert-BlankBefore "The backend server would handle the requests from the client" @($blankNoStyle)

# 2) before "The database would store the data for the application, including user accounts..."
Insert-BlankBefore "The database would store the data for the application, including user accounts" @($blankNoStyle)

# 3) before "The machine learning pipeline would be integrated into the backend server..."
Insert-BlankBefore "The machine learning pipeline would be integrated into the backend server" @($blankListStyle, $blankNoStyle)

# 4) add <w:lastRenderedPageBreak/> right before the Postprocessing run text.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Postprocessing: The process of refining"
$find.Execute() | Out-Null
if (-not $find.Found) { throw "Postprocessing run not found" }
$postPara = $find.Parent.Paragraphs(1)
$postRun = $postPara.Range
$postXml = Get-ParaXml($postRun)
$postXml2 = $postXml -replace '(<w:t[ >])', '<w:lastRenderedPageBreak/>$1', 1
$postRun.InsertXML($postXml2)

# 5) remove <w:lastRenderedPageBreak/> from the start of the Evaluation run.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "Evaluation: The process of measuring"
$find2.Execute() | Out-Null
if (-not $find2.Found) { throw "Evaluation run not found" }
$evalPara = $find2.Parent.Paragraphs(1)
$evalRange = $evalPara.Range
$evalXml = Get-ParaXml($evalRange)
$evalXml2 = $evalXml -replace '<w:lastRenderedPageBreak/>', ''
$evalRange.InsertXML($evalXml2)

Write-Output "done"
